# Weekly refresh of the "Femacal de La Calera - Pepino dulce" price table.
# Data rows (2-12) are rotated/updated in place to reflect the new weekly
# snapshot, per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 44235; I = "Primera"; J = 80; K = 14000; L = 14000; M = 14000; N = "`$/bandeja 18 kilos"; P = 778;  Q = 18 },
    @{ Row = 3;  D = 44235; I = "Segunda"; J = 70; K = 12000; L = 12000; M = 12000; N = "`$/bandeja 18 kilos"; P = 667;  Q = 18 },
    @{ Row = 4;  D = 44235; I = "Tercera"; J = 60; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; P = 556;  Q = 18 },
    @{ Row = 5;  D = 44242; I = "Primera"; J = 60; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 },
    @{ Row = 6;  D = 44242; I = "Segunda"; J = 50; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; P = 556;  Q = 18 },
    @{ Row = 7;  D = 44238; I = "Primera"; J = 90; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 },
    @{ Row = 8;  D = 44238; I = "Segunda"; J = 80; K = 11000; L = 11000; M = 11000; N = "`$/bandeja 18 kilos"; P = 611;  Q = 18 },
    @{ Row = 9;  D = 44536; I = "Primera"; J = 87; K = 22000; L = 22000; M = 22000; N = "`$/bandeja 18 kilos"; P = 1222; Q = 18 },
    @{ Row = 10; D = 44536; I = "Segunda"; J = 80; K = 20000; L = 20000; M = 20000; N = "`$/bandeja 18 kilos"; P = 1111; Q = 18 },
    @{ Row = 11; D = 44424; I = "Primera"; J = 75; K = 18000; L = 18000; M = 18000; N = "`$/caja 15 kilos";    P = 1200; Q = 15 },
    @{ Row = 12; D = 44424; I = "Segunda"; J = 50; K = 12000; L = 12000; M = 12000; N = "`$/caja 15 kilos";    P = 800;  Q = 15 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha (serial date, format preserved)
    $ws.Cells.Item($row, 9).Value  = $r.I   # I: Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Unidad de comercialización
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Kg o Unidades
}
